$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.900.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.657.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.897.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.651.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.600"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +12.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.937.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.02%  "
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.441.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("E35").Value = "  +4.48%  "
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("E37").Value = "  +2.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "78.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.92%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.845"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0501"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "50.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.803.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "94.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0108"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.30%  "
